# participant register logic added
# Minor table-grid width recalculation: the first table's first column
# (and the corresponding cell widths) shrink by 1 twip (1445 -> 1444),
# and the second table's second column (and its cell widths) shrink by
# 1 twip (7899 -> 7898). Widths in the Word object model are expressed
# in points, where 1 point = 20 twips.

$d = $word.ActiveDocument

# --- Table 1: "Git Repo" / "Branch" table ---
$t1 = $d.Tables.Item(1)
$t1.Columns.Item(1).Width = 1444 / 20

# --- Table 2: "Database details" table ---
$t2 = $d.Tables.Item(2)
$t2.Columns.Item(2).Width = 7898 / 20
